# Auto-generated edit script: updates cached market-price / profit values
# across the Cactuar_Profits workbook sheets, per scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 8181.125
$ws.Range("I70").Value = 5724.5
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 17173.5
$ws.Range("L70").Value = 27000
$ws.Range("M70").Value = -16903.5
$ws.Range("N70").Value = -27540
# Row 73
$ws.Range("H73").Value = 8181.125
$ws.Range("I73").Value = 5724.5
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 17173.5
$ws.Range("L73").Value = 27000
$ws.Range("M73").Value = -16237.5
$ws.Range("N73").Value = -28872
# Row 80
$ws.Range("H80").Value = 15625832
$ws.Range("J80").Value = 1256.8572
$ws.Range("L80").Value = 3770.5716
$ws.Range("N80").Value = -5766.571599999999
# Row 83
$ws.Range("H83").Value = 15625832
$ws.Range("J83").Value = 1256.8572
$ws.Range("L83").Value = 11311.7148
$ws.Range("N83").Value = -21295.7148
# Row 113
$ws.Range("H113").Value = 102995.14
$ws.Range("I113").Value = 4000
$ws.Range("K113").Value = 4000
$ws.Range("M113").Value = -746
# Row 125
$ws.Range("H125").Value = 4750
$ws.Range("J125").Value = 4750
$ws.Range("L125").Value = 42750
$ws.Range("N125").Value = -47670

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2704.2693
$ws.Range("J45").Value = 2872.0952
$ws.Range("L45").Value = 2872.0952
$ws.Range("N45").Value = -3626.0952
# Row 97
$ws.Range("H97").Value = 476.13333
$ws.Range("I97").Value = 251.63637
$ws.Range("K97").Value = 251.63637
$ws.Range("M97").Value = 244.36363
# Row 122
$ws.Range("H122").Value = 3475993.8
$ws.Range("I122").Value = 4764875
$ws.Range("K122").Value = 14294625
$ws.Range("M122").Value = -14292175

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 142861890
$ws.Range("I86").Value = 3802.25
$ws.Range("J86").Value = 333339330
$ws.Range("K86").Value = 3802.25
$ws.Range("L86").Value = 333339330
$ws.Range("M86").Value = -2679.25
$ws.Range("N86").Value = -333341576
# Row 89
$ws.Range("H89").Value = 142861890
$ws.Range("I89").Value = 3802.25
$ws.Range("J89").Value = 333339330
$ws.Range("K89").Value = 19011.25
$ws.Range("L89").Value = 1666696650
$ws.Range("M89").Value = -13395.25
$ws.Range("N89").Value = -1666707882

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1238.4
$ws.Range("I16").Value = 899.8
$ws.Range("K16").Value = 899.8
$ws.Range("M16").Value = -612.8
# Row 62
$ws.Range("H62").Value = 94501
$ws.Range("I62").Value = 54002
$ws.Range("K62").Value = 54002
$ws.Range("M62").Value = -53378
# Row 65
$ws.Range("H65").Value = 94501
$ws.Range("I65").Value = 54002
$ws.Range("K65").Value = 270010
$ws.Range("M65").Value = -266890
# Row 97
$ws.Range("H97").Value = 23795
$ws.Range("J97").Value = 23795
$ws.Range("L97").Value = 23795
$ws.Range("N97").Value = -25777
# Row 113
$ws.Range("H113").Value = 1238.4
$ws.Range("I113").Value = 899.8
$ws.Range("K113").Value = 899.8
$ws.Range("M113").Value = 1270.2

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 595.4583
$ws.Range("J12").Value = 775.25
$ws.Range("L12").Value = 2325.75
$ws.Range("N12").Value = -2671.75
# Row 17
$ws.Range("H17").Value = 487.5
$ws.Range("J17").Value = 875
$ws.Range("L17").Value = 2625
$ws.Range("N17").Value = -2963
# Row 39
$ws.Range("H39").Value = 4712.1816
$ws.Range("J39").Value = 4688.9
$ws.Range("L39").Value = 14066.7
$ws.Range("N39").Value = -14654.7
# Row 58
$ws.Range("H58").Value = 10666.667
$ws.Range("I58").Value = 9000
$ws.Range("J58").Value = 11500
$ws.Range("K58").Value = 27000
$ws.Range("L58").Value = 34500
$ws.Range("N58").Value = -34756
$ws.Range("M58").Value = -26872

$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 44900
$ws.Range("I58").Value = 39800
$ws.Range("K58").Value = 39800
$ws.Range("M58").Value = -39523
# Row 70
$ws.Range("H70").Value = 1901161.1
$ws.Range("I70").Value = 3253505.8
$ws.Range("K70").Value = 3253505.8
$ws.Range("M70").Value = -3253235.8
# Row 73
$ws.Range("H73").Value = 1901161.1
$ws.Range("I73").Value = 3253505.8
$ws.Range("K73").Value = 3253505.8
$ws.Range("M73").Value = -3252569.8
# Row 80
$ws.Range("H80").Value = 2376142.5
$ws.Range("J80").Value = 2166.3333
$ws.Range("L80").Value = 2166.3333
$ws.Range("N80").Value = -4162.3333
# Row 83
$ws.Range("H83").Value = 2376142.5
$ws.Range("J83").Value = 2166.3333
$ws.Range("L83").Value = 10831.6665
$ws.Range("N83").Value = -20815.6665
# Row 122
$ws.Range("H122").Value = 460120
$ws.Range("I122").Value = 1113883.6
$ws.Range("K122").Value = 3341650.8
$ws.Range("M122").Value = -3339200.8
# Row 126
$ws.Range("H126").Value = 5770.3335
$ws.Range("I126").Value = 2088.3333
$ws.Range("J126").Value = 6997.6665
$ws.Range("K126").Value = 6264.999899999999
$ws.Range("L126").Value = 20992.9995
$ws.Range("M126").Value = -3794.999899999999
$ws.Range("N126").Value = -25932.9995
# Row 132
$ws.Range("H132").Value = 4680.3335
$ws.Range("I132").Value = 4376.3125
$ws.Range("J132").Value = 5288.375
$ws.Range("K132").Value = 13128.9375
$ws.Range("L132").Value = 15865.125
$ws.Range("M132").Value = -10598.9375
$ws.Range("N132").Value = -20925.125

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3912.7307
$ws.Range("I7").Value = 3215.6604
$ws.Range("J7").Value = 5390.52
$ws.Range("K7").Value = 3215.6604
$ws.Range("L7").Value = 5390.52
$ws.Range("M7").Value = -3103.6604
$ws.Range("N7").Value = -5614.52
# Row 16
$ws.Range("H16").Value = 1568.4166
$ws.Range("J16").Value = 2715.5
$ws.Range("L16").Value = 2715.5
$ws.Range("N16").Value = -3055.5
# Row 61
$ws.Range("H61").Value = 6586.5
$ws.Range("I61").Value = 7499.6665
$ws.Range("J61").Value = 3847
$ws.Range("K61").Value = 7499.6665
$ws.Range("L61").Value = 3847
$ws.Range("M61").Value = -7297.6665
$ws.Range("N61").Value = -4251
# Row 68
$ws.Range("H68").Value = 878043.4
$ws.Range("I68").Value = 1036233.8
$ws.Range("J68").Value = 7996
$ws.Range("K68").Value = 1036233.8
$ws.Range("L68").Value = 7996
$ws.Range("M68").Value = -1035484.8
$ws.Range("N68").Value = -9494
# Row 71
$ws.Range("H71").Value = 878043.4
$ws.Range("I71").Value = 1036233.8
$ws.Range("J71").Value = 7996
$ws.Range("K71").Value = 5181169
$ws.Range("L71").Value = 39980
$ws.Range("M71").Value = -5177425
$ws.Range("N71").Value = -47468
# Row 113
$ws.Range("H113").Value = 6586.5
$ws.Range("I113").Value = 7499.6665
$ws.Range("J113").Value = 3847
$ws.Range("K113").Value = 7499.6665
$ws.Range("L113").Value = 3847
$ws.Range("M113").Value = -5329.6665
$ws.Range("N113").Value = -8187
# Row 122
$ws.Range("H122").Value = 7670.476
$ws.Range("I122").Value = 2548.077
$ws.Range("J122").Value = 15994.375
$ws.Range("K122").Value = 7644.231000000001
$ws.Range("L122").Value = 47983.125
$ws.Range("M122").Value = -5194.231000000001
$ws.Range("N122").Value = -52883.125
# Row 126
$ws.Range("H126").Value = 3912.7307
$ws.Range("I126").Value = 3215.6604
$ws.Range("J126").Value = 5390.52
$ws.Range("K126").Value = 9646.9812
$ws.Range("L126").Value = 16171.56
$ws.Range("M126").Value = -7176.9812
$ws.Range("N126").Value = -21111.56
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 12700
$ws.Range("J62").Value = 12700
$ws.Range("L62").Value = 12700
$ws.Range("N62").Value = -13948
# Row 65
$ws.Range("H65").Value = 12700
$ws.Range("J65").Value = 12700
$ws.Range("L65").Value = 63500
$ws.Range("N65").Value = -69740

Write-Host "Done: updated cells across 8 sheets"